# Scaling bed time data to am/pm
# The "Bed Time" column (O) contains free-text time strings (HH:MM:SS, 24h style).
# This scales a handful of 24-hour-looking entries down to their 12-hour
# equivalents (e.g. 22:00:00 -> 10:00:00), per the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value   = "00:00:00"
$ws.Range("O8").Value   = "10:00:00"
$ws.Range("O9").Value   = "00:30:00"
$ws.Range("O10").Value  = "11:30:00"
$ws.Range("O12").Value  = "11:00:00"
$ws.Range("O16").Value  = "00:30:00"
$ws.Range("O23").Value  = "00:30:00"
$ws.Range("O24").Value  = "00:30:00"
$ws.Range("O25").Value  = "09:00:00"
$ws.Range("O26").Value  = "11:00:00"
$ws.Range("O35").Value  = "00:30:00"
$ws.Range("O38").Value  = "11:30:00"
$ws.Range("O40").Value  = "11:00:00"
$ws.Range("O41").Value  = "11:00:00"
$ws.Range("O56").Value  = "00:00:00"
$ws.Range("O60").Value  = "00:30:00"
$ws.Range("O63").Value  = "11:00:00"
$ws.Range("O70").Value  = "00:00:00"
$ws.Range("O81").Value  = "11:30:00"
$ws.Range("O82").Value  = "11:00:00"
$ws.Range("O83").Value  = "11:00:00"
$ws.Range("O92").Value  = "11:30:00"
$ws.Range("O94").Value  = "10:00:00"
$ws.Range("O95").Value  = "11:00:00"
$ws.Range("O104").Value = "09:00:00"
